$wb = $excel.ActiveWorkbook

$appsV2 = $wb.Worksheets.Item("Apps-SMB-V2")
$torture = $wb.Worksheets.Item("SMB2TORTURE")

# ---------------------------------------------------------------------------
# Style source cells (already-existing cells carrying the exact named
# styles we need), taken from the Apps-SMB-V2 legend block (rows 25-28)
# and the matrix body (row 2) so every new cell reuses an existing cellXf
# instead of Excel minting a brand-new one.
# ---------------------------------------------------------------------------
# B25 -> "Done"            (style 12, glyph referencing shared string 314 "ü")
# B26 -> "Partially Done"  (style 14, glyph referencing shared string 317 "?")
# B27 -> "Pending"         (style 13, glyph referencing shared string 315 "û")
# B28 -> "Neutral"         (style 15, glyph referencing shared string 316 "!")
# C25..C28 -> "Normal", centered (style 7)

# ---------------------------------------------------------------------------
# SMB2TORTURE (sheet10): fill in column C ("N/A" -> mostly "Partially Done")
# for rows 2-21, with row 10 marked "Done".
# ---------------------------------------------------------------------------
$partialRows = @(2,3,4,5,6,7,8,9,11,12,13,14,15,16,17,18,19,20,21)
foreach ($r in $partialRows) {
    $appsV2.Range("B26").Copy()
    $torture.Range("C$r").PasteSpecial(-4122)
    $torture.Range("C$r").Value = "?"
}

$appsV2.Range("B25").Copy()
$torture.Range("C10").PasteSpecial(-4122)
$torture.Range("C10").Value = "ü"

# Legend block, rows 23-26, columns B (glyph) / C (description)
$appsV2.Range("B25").Copy()
$torture.Range("B23").PasteSpecial(-4122)
$torture.Range("B23").Value = "ü"
$appsV2.Range("C25").Copy()
$torture.Range("C23").PasteSpecial(-4122)
$torture.Range("C23").Value = "Works"

$appsV2.Range("B26").Copy()
$torture.Range("B24").PasteSpecial(-4122)
$torture.Range("B24").Value = "?"
$appsV2.Range("C26").Copy()
$torture.Range("C24").PasteSpecial(-4122)
$torture.Range("C24").Value = "T.B.D."

$appsV2.Range("B27").Copy()
$torture.Range("B25").PasteSpecial(-4122)
$torture.Range("B25").Value = "û"
$appsV2.Range("C27").Copy()
$torture.Range("C25").PasteSpecial(-4122)
$torture.Range("C25").Value = "Fails"

$appsV2.Range("B28").Copy()
$torture.Range("B26").PasteSpecial(-4122)
$torture.Range("B26").Value = [char]0x00A1
$appsV2.Range("C28").Copy()
$torture.Range("C26").PasteSpecial(-4122)
$torture.Range("C26").Value = "N/A"

# Column D on SMB2TORTURE widens from being merged with column C
# (min=3 max=4) to its own explicit width (~21.57 characters).
$torture.Columns.Item(4).ColumnWidth = 20.65

# ---------------------------------------------------------------------------
# Apps-SMB-V2 (sheet9): C18 flips from "Partially Done" to "Pending".
# ---------------------------------------------------------------------------
$appsV2.Range("B27").Copy()
$appsV2.Range("C18").PasteSpecial(-4122)
$appsV2.Range("C18").Value = "û"

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping to match the recorded UI state.
# ---------------------------------------------------------------------------
$appsV2.Activate()
$appsV2.Range("C18").Select()

$torture.Activate()
$torture.Range("C21").Select()
